# Reorder the "Recorded By" (column G) values on the "Session Analysis Results"
# sheet: for every row whose G cell contains multiple comma-separated
# recorder names, reverse the order of the list (e.g. "System, foo@bar.com"
# becomes "foo@bar.com, System"). Rows with a single recorder are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $current = $cell.Value2

    if ($current -ne $null -and $current -like "*, *") {
        $parts = $current -split ", "
        $reversed = $parts[($parts.Length - 1)..0]
        $newValue = $reversed -join ", "
        $cell.Value = $newValue
    }
}
